$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1527.3115
$ws.Range("I132").Value = 1359.4
$ws.Range("K132").Value = 4078.2
$ws.Range("M132").Value = -1548.2
$ws.Range("H138").Value = 3775.6567
$ws.Range("I138").Value = 1523.3846
$ws.Range("J138").Value = 4317.8706
$ws.Range("K138").Value = 4570.1538
$ws.Range("L138").Value = 12953.6118
$ws.Range("M138").Value = 569.8462
$ws.Range("N138").Value = -23233.6118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 20000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 20000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 20000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -20700
$ws.Range("H17").Value = 100009
$ws.Range("J17").Value = 100009
$ws.Range("L17").Value = 100009
$ws.Range("N17").Value = -100355
$ws.Range("H41").Value = 5871.2
$ws.Range("I41").Value = 1139
$ws.Range("J41").Value = 24800
$ws.Range("K41").Value = 1139
$ws.Range("L41").Value = 24800
$ws.Range("M41").Value = -725
$ws.Range("N41").Value = -25628
$ws.Range("H110").Value = 1166.6666
$ws.Range("I110").Value = 500
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 500
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 1545
$ws.Range("N110").Value = -5590
$ws.Range("H122").Value = 4312433
$ws.Range("I122").Value = 2608.8572
$ws.Range("J122").Value = 5683741
$ws.Range("K122").Value = 7826.571599999999
$ws.Range("L122").Value = 17051223
$ws.Range("M122").Value = -5376.571599999999
$ws.Range("N122").Value = -17056123
$ws.Range("H132").Value = 5357
$ws.Range("I132").Value = 1237.0555
$ws.Range("K132").Value = 3711.1665
$ws.Range("M132").Value = -1181.1665
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 10844.857
$ws.Range("I75").Value = 7319
$ws.Range("K75").Value = 7319
$ws.Range("M75").Value = -6383
$ws.Range("H78").Value = 10844.857
$ws.Range("I78").Value = 7319
$ws.Range("K78").Value = 21957
$ws.Range("M78").Value = -17277
$ws.Range("H86").Value = 2391.5652
$ws.Range("I86").Value = 2391.5652
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2391.5652
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1268.5652
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2391.5652
$ws.Range("I89").Value = 2391.5652
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 11957.826
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6341.826000000001
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 7866.909
$ws.Range("I134").Value = 7753.6
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 23260.8
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -20725.8
$ws.Range("N134").Value = -32070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3124.3555
$ws.Range("I31").Value = 2590.0908
$ws.Range("J31").Value = 3635.3914
$ws.Range("K31").Value = 2590.0908
$ws.Range("L31").Value = 3635.3914
$ws.Range("M31").Value = -2295.0908
$ws.Range("N31").Value = -4225.3914
$ws.Range("H34").Value = 3124.3555
$ws.Range("I34").Value = 2590.0908
$ws.Range("J34").Value = 3635.3914
$ws.Range("K34").Value = 2590.0908
$ws.Range("L34").Value = 3635.3914
$ws.Range("M34").Value = -2388.0908
$ws.Range("N34").Value = -4039.3914
$ws.Range("H58").Value = 2023328.6
$ws.Range("I58").Value = 3137132.2
$ws.Range("J58").Value = 4559.375
$ws.Range("K58").Value = 3137132.2
$ws.Range("L58").Value = 4559.375
$ws.Range("M58").Value = -3136929.2
$ws.Range("N58").Value = -4965.375
$ws.Range("H136").Value = 2023328.6
$ws.Range("I136").Value = 3137132.2
$ws.Range("J136").Value = 4559.375
$ws.Range("K136").Value = 9411396.600000001
$ws.Range("L136").Value = 13678.125
$ws.Range("M136").Value = -9408846.600000001
$ws.Range("N136").Value = -18778.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 666.8333
$ws.Range("I11").Value = 275.25
$ws.Range("J11").Value = 1450
$ws.Range("K11").Value = 825.75
$ws.Range("L11").Value = 4350
$ws.Range("M11").Value = -685.75
$ws.Range("N11").Value = -4630
$ws.Range("H114").Value = 458.7
$ws.Range("I114").Value = 278.33334
$ws.Range("J114").Value = 999.8
$ws.Range("K114").Value = 835.0000200000001
$ws.Range("L114").Value = 2999.4
$ws.Range("M114").Value = 2418.99998
$ws.Range("N114").Value = -9507.4
$ws.Range("H131").Value = 34721.82
$ws.Range("I131").Value = 647.7059
$ws.Range("J131").Value = 87381.82000000001
$ws.Range("K131").Value = 1943.1177
$ws.Range("L131").Value = 262145.46
$ws.Range("M131").Value = 3096.8823
$ws.Range("N131").Value = -272225.46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4759.8955
$ws.Range("I70").Value = 4251
$ws.Range("J70").Value = 4992.2173
$ws.Range("K70").Value = 4251
$ws.Range("L70").Value = 4992.2173
$ws.Range("M70").Value = -3981
$ws.Range("N70").Value = -5532.2173
$ws.Range("H73").Value = 4759.8955
$ws.Range("I73").Value = 4251
$ws.Range("J73").Value = 4992.2173
$ws.Range("K73").Value = 4251
$ws.Range("L73").Value = 4992.2173
$ws.Range("M73").Value = -3315
$ws.Range("N73").Value = -6864.2173
$ws.Range("H102").Value = 2977.8372
$ws.Range("I102").Value = 2707.4517
$ws.Range("K102").Value = 2707.4517
$ws.Range("M102").Value = -1085.4517
$ws.Range("H132").Value = 6790.4287
$ws.Range("I132").Value = 9048.385
$ws.Range("K132").Value = 27145.155
$ws.Range("M132").Value = -24615.155
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 35433.332
$ws.Range("J64").Value = 35433.332
$ws.Range("L64").Value = 35433.332
$ws.Range("N64").Value = -35883.332
$ws.Range("H67").Value = 35433.332
$ws.Range("J67").Value = 35433.332
$ws.Range("L67").Value = 35433.332
$ws.Range("N67").Value = -36993.332
$ws.Range("H122").Value = 5182.3096
$ws.Range("I122").Value = 4544.731
$ws.Range("J122").Value = 6218.375
$ws.Range("K122").Value = 13634.193
$ws.Range("L122").Value = 18655.125
$ws.Range("M122").Value = -11184.193
$ws.Range("N122").Value = -23555.125
$ws.Range("H132").Value = 2376.68
$ws.Range("I132").Value = 1735.1538
$ws.Range("J132").Value = 3071.6667
$ws.Range("K132").Value = 5205.4614
$ws.Range("L132").Value = 9215.000100000001
$ws.Range("M132").Value = -2675.4614
$ws.Range("N132").Value = -14275.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 12328.454
$ws.Range("J54").Value = 12328.454
$ws.Range("L54").Value = 12328.454
$ws.Range("N54").Value = -13368.454
$ws.Range("H63").Value = 40249
$ws.Range("J63").Value = 40249
$ws.Range("L63").Value = 40249
$ws.Range("N63").Value = -41497
$ws.Range("H66").Value = 40249
$ws.Range("J66").Value = 40249
$ws.Range("L66").Value = 120747
$ws.Range("N66").Value = -126987
$ws.Range("H132").Value = 1432.2703
$ws.Range("I132").Value = 1236.3667
$ws.Range("J132").Value = 2271.8572
$ws.Range("K132").Value = 3709.1001
$ws.Range("L132").Value = 6815.571599999999
$ws.Range("M132").Value = -1179.1001
$ws.Range("N132").Value = -11875.5716
$ws.Range("H136").Value = 6186.5854
$ws.Range("I136").Value = 2660.05
$ws.Range("J136").Value = 9545.190000000001
$ws.Range("K136").Value = 7980.150000000001
$ws.Range("L136").Value = 28635.57
$ws.Range("M136").Value = -5430.150000000001
$ws.Range("N136").Value = -33735.57
